$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# queueAnimation / sticky rows no longer link to a CodePen example -> "N/A"
$ws.Range("D15").Value = "N/A"
$ws.Range("D16").Value = "N/A"

# Point the phase-options descriptions at the corrected anchor in the README
$ws.Range("E12").Value = "Accepts an Object that will pass options whenever you animate to the start phase. For a list of valid phase options refer to the [Options (Phase)](https://github.com/RealTayy/slice-revealer#options---phase-optional) section."
$ws.Range("E13").Value = "Accepts an Object that will pass options whenever you animate to the halfway phase. For a list of valid phase options refer to the [Options (Phase)](https://github.com/RealTayy/slice-revealer#options---phase-optional) section."
$ws.Range("E14").Value = "Accepts an Object that will pass options whenever you animate to the end phase. For a list of valid phase options refer to the [Options (Phase)](https://github.com/RealTayy/slice-revealer#options---phase-optional) section."

# Fix "Slive Revealer" typo -> "Slice Revealer"
$ws.Range("E16").Value = "By default the position of the element that contains the Slice Revealer instance has ``position: abosolute``. If true then container element will have ``position: sticky``. Used when target has scrollable content."

# Leave the selection where the author last clicked while editing
$ws.Range("G20").Select()
